$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 17:47"
$ws.Range("B4").Value = 7011829
$ws.Range("C4").Value = 7061
$ws.Range("D4").Value = 4261777
$ws.Range("E4").Value = 2545871
$ws.Range("G4").Value = 63
$ws.Range("H4").Value = 204181
$ws.Range("B5").Value = 5523917
$ws.Range("C5").Value = 38305
$ws.Range("D5").Value = 4440775
$ws.Range("E5").Value = 994797
$ws.Range("G5").Value = 436
$ws.Range("H5").Value = 88345
$ws.Range("B15").Value = 447468
$ws.Range("C15").Value = 1194
$ws.Range("D15").Value = 421111
$ws.Range("E15").Value = 14059
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 12298
$ws.Range("B17").Value = 398625
$ws.Range("C17").Value = 4368
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 41788
$ws.Range("B23").Value = 299506
$ws.Range("C23").Value = 1350
$ws.Range("D23").Value = 218703
$ws.Range("E23").Value = 45079
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = 35724
$ws.Range("B25").Value = 274717
$ws.Range("C25").Value = 1240
$ws.Range("E25").Value = 21241
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 9476
$ws.Range("B29").Value = 144076
$ws.Range("C29").Value = 427
$ws.Range("D29").Value = 124869
$ws.Range("E29").Value = 9988
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 9219
$ws.Range("B34").Value = 108783
$ws.Range("C34").Value = 494
$ws.Range("D34").Value = 82274
$ws.Range("E34").Value = 24455
$ws.Range("G34").Value = 7
$ws.Range("H34").Value = 2054
$ws.Range("A44").Value = "Guatemala"
$ws.Range("B44").Value = 85681
$ws.Range("C44").Value = 237
$ws.Range("D44").Value = 75172
$ws.Range("E44").Value = 7385
$ws.Range("G44").Value = 5
$ws.Range("H44").Value = 3124
$ws.Range("A45").Value = "Emiratos Arabes Unidos"
$ws.Range("B45").Value = 85595
$ws.Range("C45").Value = 679
$ws.Range("D45").Value = 75086
$ws.Range("E45").Value = 10104
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 405
$ws.Range("B67").Value = 39280
$ws.Range("C67").Value = 92
$ws.Range("D67").Value = 36836
$ws.Range("E67").Value = 1868
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 576
$ws.Range("B70").Value = 37079
$ws.Range("C70").Value = 98
$ws.Range("E70").Value = 12544
$ws.Range("B82").Value = 20598
$ws.Range("C82").Value = 167
$ws.Range("E82").Value = 1058
$ws.Range("B94").Value = 12535
$ws.Range("C94").Value = 150
$ws.Range("D94").Value = 6995
$ws.Range("E94").Value = 5176
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 364
$ws.Range("B101").Value = 9712
$ws.Range("C101").Value = 20
$ws.Range("D101").Value = 9373
$ws.Range("E101").Value = 274
$ws.Range("A104").Value = "Montenegro"
$ws.Range("B104").Value = 8842
$ws.Range("C104").Value = 230
$ws.Range("D104").Value = 5425
$ws.Range("E104").Value = 3279
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 138
$ws.Range("A105").Value = "Gabon"
$ws.Range("B105").Value = 8696
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 7848
$ws.Range("E105").Value = 795
$ws.Range("H105").Value = 53
$ws.Range("A106").Value = "Haiti"
$ws.Range("B106").Value = 8624
$ws.Range("C106").Value = 5
$ws.Range("D106").Value = 6482
$ws.Range("E106").Value = 1921
$ws.Range("H106").Value = 221
$ws.Range("B107").Value = 7916
$ws.Range("C107").Value = 9
$ws.Range("D107").Value = 6839
$ws.Range("E107").Value = 953
$ws.Range("B113").Value = 6151
$ws.Range("C113").Value = 610
$ws.Range("D113").Value = 1445
$ws.Range("E113").Value = 4608
$ws.Range("G113").Value = 6
$ws.Range("H113").Value = 98
$ws.Range("B114").Value = 5733
$ws.Range("C114").Value = 2
$ws.Range("D114").Value = 4053
$ws.Range("E114").Value = 1501
$ws.Range("B115").Value = 5404
$ws.Range("C115").Value = 1
$ws.Range("D115").Value = 5336
$ws.Range("A118").Value = "Jamaica"
$ws.Range("B118").Value = 5143
$ws.Range("C118").Value = 155
$ws.Range("D118").Value = 1407
$ws.Range("E118").Value = 3666
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 70
$ws.Range("A119").Value = "Cuba"
$ws.Range("B119").Value = 5141
$ws.Range("C119").Value = 50
$ws.Range("D119").Value = 4462
$ws.Range("E119").Value = 563
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 116
$ws.Range("A120").Value = "Hong Kong"
$ws.Range("B120").Value = 5039
$ws.Range("C120").Value = 6
$ws.Range("D120").Value = 4717
$ws.Range("E120").Value = 219
$ws.Range("H120").Value = 103
$ws.Range("A121").Value = "Congo"
$ws.Range("B121").Value = 5002
$ws.Range("C121").Value = 16
$ws.Range("D121").Value = 3887
$ws.Range("E121").Value = 1026
$ws.Range("H121").Value = 89
$ws.Range("A122").Value = "Guinea Ecuatorial"
$ws.Range("B122").Value = 5002
$ws.Range("D122").Value = 4509
$ws.Range("E122").Value = 410
$ws.Range("H122").Value = 83
$ws.Range("B130").Value = 3930
$ws.Range("C130").Value = 29
$ws.Range("D130").Value = 1802
$ws.Range("E130").Value = 2063
$ws.Range("B159").Value = 1603
$ws.Range("C159").Value = 3
$ws.Range("E159").Value = 212
$ws.Range("B183").Value = 340
$ws.Range("C183").Value = 1
$ws.Range("E183").Value = 4
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
